# Auto-generated edit script applying the scheduled-runner market-price data refresh
# described by the OOXML diff (Sheets/Phoenix_Profits.xlsx).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 43890.957
$ws.Range("I98").Value = 52163.473
$ws.Range("J98").Value = 4596.5
$ws.Range("K98").Value = 52163.473
$ws.Range("L98").Value = 4596.5
$ws.Range("M98").Value = -50665.473
$ws.Range("N98").Value = -7592.5
$ws.Range("H112").Value = 1950.6666
$ws.Range("J112").Value = 2000.8889
$ws.Range("L112").Value = 6002.6667
$ws.Range("N112").Value = -8218.6667
$ws.Range("H122").Value = 43890.957
$ws.Range("I122").Value = 52163.473
$ws.Range("J122").Value = 4596.5
$ws.Range("K122").Value = 156490.419
$ws.Range("L122").Value = 13789.5
$ws.Range("M122").Value = -154040.419
$ws.Range("N122").Value = -18689.5
$ws.Range("H132").Value = 3716.1428
$ws.Range("I132").Value = 3904.7368
$ws.Range("J132").Value = 1924.5
$ws.Range("K132").Value = 11714.2104
$ws.Range("L132").Value = 5773.5
$ws.Range("M132").Value = -9184.2104
$ws.Range("N132").Value = -10833.5
$ws.Range("H137").Value = 173050.77
$ws.Range("I137").Value = 985.3226
$ws.Range("K137").Value = 2955.9678
$ws.Range("M137").Value = -405.9677999999999
$ws.Range("H138").Value = 3002.6086
$ws.Range("J138").Value = 3531.2068
$ws.Range("L138").Value = 10593.6204
$ws.Range("N138").Value = -20873.6204

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1269.6586
$ws.Range("I2").Value = 1160.8286
$ws.Range("K2").Value = 1160.8286
$ws.Range("M2").Value = -1047.8286
$ws.Range("H32").Value = 3260.7144
$ws.Range("I32").Value = 2985.1875
$ws.Range("K32").Value = 2985.1875
$ws.Range("M32").Value = -2698.1875
$ws.Range("H61").Value = 3327.5774
$ws.Range("I61").Value = 2436.4888
$ws.Range("K61").Value = 2436.4888
$ws.Range("M61").Value = -2224.4888
$ws.Range("H116").Value = 1269.6586
$ws.Range("I116").Value = 1160.8286
$ws.Range("K116").Value = 1160.8286
$ws.Range("M116").Value = 1133.1714
$ws.Range("H132").Value = 2243.8206
$ws.Range("I132").Value = 2185
$ws.Range("K132").Value = 6555
$ws.Range("M132").Value = -4025
$ws.Range("H136").Value = 3327.5774
$ws.Range("I136").Value = 2436.4888
$ws.Range("K136").Value = 7309.4664
$ws.Range("M136").Value = -4759.4664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1269.6586
$ws.Range("I3").Value = 1160.8286
$ws.Range("K3").Value = 1160.8286
$ws.Range("M3").Value = -1046.8286
$ws.Range("H35").Value = 1783.5
$ws.Range("I35").Value = 1783.5
$ws.Range("K35").Value = 1783.5
$ws.Range("M35").Value = -1473.5
$ws.Range("H105").Value = 25003472
$ws.Range("I105").Value = 35717350
$ws.Range("K105").Value = 35717350
$ws.Range("M105").Value = -35715603
$ws.Range("H134").Value = 2695.6047
$ws.Range("I134").Value = 2254.4062
$ws.Range("K134").Value = 6763.2186
$ws.Range("M134").Value = -4228.2186

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2183
$ws.Range("I31").Value = 1899.7587
$ws.Range("K31").Value = 1899.7587
$ws.Range("M31").Value = -1604.7587
$ws.Range("H34").Value = 2183
$ws.Range("I34").Value = 1899.7587
$ws.Range("K34").Value = 1899.7587
$ws.Range("M34").Value = -1697.7587
$ws.Range("H57").Value = 23417.25
$ws.Range("J57").Value = 22890.334
$ws.Range("L57").Value = 22890.334
$ws.Range("N57").Value = -24010.334
$ws.Range("H58").Value = 1922.1136
$ws.Range("I58").Value = 1505.5151
$ws.Range("K58").Value = 1505.5151
$ws.Range("M58").Value = -1302.5151
$ws.Range("H86").Value = 6559.1
$ws.Range("I86").Value = 3599
$ws.Range("J86").Value = 7827.7144
$ws.Range("K86").Value = 3599
$ws.Range("L86").Value = 7827.7144
$ws.Range("M86").Value = -2476
$ws.Range("N86").Value = -10073.7144
$ws.Range("H89").Value = 6559.1
$ws.Range("I89").Value = 3599
$ws.Range("J89").Value = 7827.7144
$ws.Range("K89").Value = 17995
$ws.Range("L89").Value = 39138.572
$ws.Range("M89").Value = -12379
$ws.Range("N89").Value = -50370.572
$ws.Range("H132").Value = 9145.071
$ws.Range("I132").Value = 3029.261
$ws.Range("J132").Value = 37277.8
$ws.Range("K132").Value = 9087.782999999999
$ws.Range("L132").Value = 111833.4
$ws.Range("M132").Value = -6557.782999999999
$ws.Range("N132").Value = -116893.4
$ws.Range("H134").Value = 3721.0476
$ws.Range("I134").Value = 4731.4
$ws.Range("J134").Value = 1195.1666
$ws.Range("K134").Value = 14194.2
$ws.Range("L134").Value = 3585.4998
$ws.Range("M134").Value = -11659.2
$ws.Range("N134").Value = -8655.4998
$ws.Range("H136").Value = 1922.1136
$ws.Range("I136").Value = 1505.5151
$ws.Range("K136").Value = 4516.5453
$ws.Range("M136").Value = -1966.5453

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 50412244
$ws.Range("I4").Value = 27889404
$ws.Range("J4").Value = 98107660
$ws.Range("K4").Value = 83668212
$ws.Range("L4").Value = 294322980
$ws.Range("M4").Value = -83668100
$ws.Range("N4").Value = -294323204
$ws.Range("H5").Value = 1166.3334
$ws.Range("I5").Value = 1166.3334
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3499.0002
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -3387.0002
$ws.Range("N5").ClearContents()
$ws.Range("H59").Value = 1269.4445
$ws.Range("I59").Value = 1178.25
$ws.Range("J59").Value = 1999
$ws.Range("K59").Value = 3534.75
$ws.Range("L59").Value = 5997
$ws.Range("M59").Value = -2994.75
$ws.Range("N59").Value = -7077
$ws.Range("H110").Value = 5798.5
$ws.Range("I110").Value = 5798.5
$ws.Range("K110").Value = 17395.5
$ws.Range("M110").Value = -13305.5
$ws.Range("H111").Value = 1508.8
$ws.Range("I111").Value = 1508.8
$ws.Range("K111").Value = 4526.4
$ws.Range("M111").Value = -1459.4
$ws.Range("H116").Value = 5915.8335
$ws.Range("I116").Value = 2998.3333
$ws.Range("J116").Value = 8833.333000000001
$ws.Range("K116").Value = 8994.999899999999
$ws.Range("L116").Value = 26499.999
$ws.Range("M116").Value = -5552.999899999999
$ws.Range("N116").Value = -33383.999
$ws.Range("H135").Value = 1166.3334
$ws.Range("I135").Value = 1166.3334
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 10497.0006
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -7962.000599999999
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 34091.28
$ws.Range("J126").Value = 3633
$ws.Range("L126").Value = 10899
$ws.Range("N126").Value = -15839
$ws.Range("H132").Value = 2816.9656
$ws.Range("I132").Value = 2997.6956
$ws.Range("J132").Value = 2124.1667
$ws.Range("K132").Value = 8993.086800000001
$ws.Range("L132").Value = 6372.500100000001
$ws.Range("M132").Value = -6463.086800000001
$ws.Range("N132").Value = -11432.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 38000
$ws.Range("J54").Value = 38000
$ws.Range("L54").Value = 38000
$ws.Range("N54").Value = -39288
$ws.Range("H82").Value = 1158.7858
$ws.Range("I82").Value = 931
$ws.Range("J82").Value = 1994
$ws.Range("K82").Value = 931
$ws.Range("L82").Value = 1994
$ws.Range("M82").Value = -570
$ws.Range("N82").Value = -2716
$ws.Range("H85").Value = 1158.7858
$ws.Range("I85").Value = 931
$ws.Range("J85").Value = 1994
$ws.Range("K85").Value = 931
$ws.Range("L85").Value = 1994
$ws.Range("M85").Value = 317
$ws.Range("N85").Value = -4490
$ws.Range("H132").Value = 3085.2163
$ws.Range("I132").Value = 2252.25
$ws.Range("K132").Value = 6756.75
$ws.Range("M132").Value = -4226.75
$ws.Range("H136").Value = 21192.479
$ws.Range("I136").Value = 2094.923
$ws.Range("J136").Value = 103948.555
$ws.Range("K136").Value = 6284.768999999999
$ws.Range("L136").Value = 311845.665
$ws.Range("M136").Value = -3734.768999999999
$ws.Range("N136").Value = -316945.665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 62449.902
$ws.Range("I136").Value = 55065.39
$ws.Range("K136").Value = 165196.17
$ws.Range("M136").Value = -162646.17
